$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.2
$ws.Range("I2").Value = 3.25
$ws.Range("J2").Value = 1.04
$ws.Range("K2").Value = 13
$ws.Range("Z2").Value = 13

# Row 4
$ws.Range("G4").Value = 1.42
$ws.Range("H4").Value = 3.95
$ws.Range("I4").Value = 7.3
$ws.Range("L4").Value = 1.28
$ws.Range("M4").Value = 3.05
$ws.Range("N4").Value = 1.83
$ws.Range("O4").Value = 1.78
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 1.65
$ws.Range("T4").Value = 6
$ws.Range("U4").Value = 6.3
$ws.Range("V4").Value = 8.25
$ws.Range("X4").Value = 12
$ws.Range("Y4").Value = 30
$ws.Range("Z4").Value = 9.5
$ws.Range("AA4").Value = 7.9
$ws.Range("AB4").Value = 21
$ws.Range("AC4").Value = 120
$ws.Range("AE4").Value = 16
$ws.Range("AJ4").Value = 90

# Row 6
$ws.Range("G6").Value = 2
$ws.Range("I6").Value = 3.75
$ws.Range("P6").Value = 1.47
$ws.Range("Q6").Value = 2.32
$ws.Range("T6").Value = 6
$ws.Range("U6").Value = 8.75
$ws.Range("W6").Value = 18
$ws.Range("Z6").Value = 7.3
$ws.Range("AE6").Value = 8.75
$ws.Range("AF6").Value = 19

# Row 7
$ws.Range("G7").Value = 2.57
$ws.Range("H7").Value = 3.4
$ws.Range("L7").Value = 1.3
$ws.Range("M7").Value = 2.92
$ws.Range("N7").Value = 1.88
$ws.Range("O7").Value = 1.72
$ws.Range("R7").Value = 1.72
$ws.Range("S7").Value = 1.88
$ws.Range("T7").Value = 8.25
$ws.Range("W7").Value = 27
$ws.Range("X7").Value = 22
$ws.Range("Y7").Value = 32
$ws.Range("Z7").Value = 10
$ws.Range("AA7").Value = 6.5
$ws.Range("AB7").Value = 14.5
$ws.Range("AC7").Value = 70
$ws.Range("AE7").Value = 8.25
$ws.Range("AF7").Value = 12
